$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old B:E columns (shift F,G into B,C) so the sheet keeps the
# "tail" headers (3932,3933) and their data, matching the uploaded file.
$ws.Range("B1:E5").Delete(-4159)

# New D/E columns with brand-new header + data values.
$ws.Range("D1").Value = 3934
$ws.Range("E1").Value = 3935

$ws.Range("D2").Value = 31.7
$ws.Range("E2").Value = 38.1

$ws.Range("D3").Value = 29.4
$ws.Range("E3").Value = 37.2

$ws.Range("D4").Value = 31.5
$ws.Range("E4").Value = 38.2

$ws.Range("D5").Value = 30.9
$ws.Range("E5").Value = 38.2

# Apply a thin box border around every cell of the table first ...
$table = $ws.Range("B1:E5")
$table.Borders.LineStyle = 1
$table.Borders.Color = 0

# ... then thicken the header row's outer edge (top/bottom/left/right).
$header = $ws.Range("B1:E1")
$header.Borders.Item(8).Weight = -4138
$header.Borders.Item(8).Color = 0
$header.Borders.Item(9).Weight = -4138
$header.Borders.Item(9).Color = 0
$ws.Range("B1").Borders.Item(7).Weight = -4138
$ws.Range("B1").Borders.Item(7).Color = 0
$ws.Range("E1").Borders.Item(10).Weight = -4138
$ws.Range("E1").Borders.Item(10).Color = 0

# The first body row sits flush against the header's bottom border, so its
# own top edge is cleared.
$ws.Range("B2:E2").Borders.Item(8).LineStyle = -4142

$ws.Range("D11").Select()
